{"js": "// Update the Mantel-test \"p\" column values in the correlogram table to\n// reflect the regenerated (smaller) correlograms.\nconst replacements = [\n  [\"0.221\", \"0.203\"],\n  [\"0.192\", \"0.172\"],\n  [\"0.494\", \"0.467\"],\n  [\"0.662\", \"0.608\"],\n  [\"0.543\", \"0.476\"],\n  [\"0.609\", \"0.515\"],\n  [\"0.883\", \"0.811\"],\n  [\"0.18\", \"0.17\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the Mantel-test \"p\" column values in the correlogram table to\n# reflect the regenerated (smaller) correlograms.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"0.221\"; New = \"0.203\" },\n    @{ Old = \"0.192\"; New = \"0.172\" },\n    @{ Old = \"0.494\"; New = \"0.467\" },\n    @{ Old = \"0.662\"; New = \"0.608\" },\n    @{ Old = \"0.543\"; New = \"0.476\" },\n    @{ Old = \"0.609\"; New = \"0.515\" },\n    @{ Old = \"0.883\"; New = \"0.811\" },\n    @{ Old = \"0.18\";  New = \"0.17\" }\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    [void]$range.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $r.New, $wdReplaceAll)\n}\n"}
